# Update "想去人数" (column F) figures on the 展览 and 全部类型 sheets to
# match the newly scraped numbers.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 7889
$ws1.Range("F5").Value  = 7889
$ws1.Range("F9").Value  = 8648
$ws1.Range("F10").Value = 8648
$ws1.Range("F13").Value = 92
$ws1.Range("F16").Value = 2778
$ws1.Range("F21").Value = 623
$ws1.Range("F22").Value = 93
$ws1.Range("F23").Value = 3952
$ws1.Range("F30").Value = 5593
$ws1.Range("F37").Value = 2825
$ws1.Range("F40").Value = 1125
$ws1.Range("F45").Value = 3635
$ws1.Range("F47").Value = 2349

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 7889
$ws4.Range("F6").Value  = 7889
$ws4.Range("F9").Value  = 8648
$ws4.Range("F10").Value = 8648
$ws4.Range("F12").Value = 92
$ws4.Range("F15").Value = 2778
$ws4.Range("F21").Value = 623
$ws4.Range("F22").Value = 93
$ws4.Range("F23").Value = 3952
$ws4.Range("F30").Value = 5593
$ws4.Range("F37").Value = 2826
$ws4.Range("F40").Value = 1125
$ws4.Range("F46").Value = 3635
$ws4.Range("F47").Value = 2349
